# Update 1.6.5.6: Bug fix
# Clear the hard-coded driver name placeholder out of C6 (it was left over
# from testing and should be blank so the template stays generic).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").ClearContents()

# Row 6 no longer needs the extra height that was reserved for the name.
$ws.Rows("6").RowHeight = 15

# Row 8 goes back to the sheet's default (auto) height.
$ws.Rows("8").AutoFit()

# Move the active selection back to the now-empty cell.
[void]$ws.Range("C6").Select()
